$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K1").Value = "SMCode"
$ws.Range("K1").Select()
